$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the team record columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style already used by the
# other header cells in row 1 (e.g. A1..AC1) by copying A1's formatting.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Team record values for every data row (2-44): 92 wins, 70 losses, 0 ties
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 92
    $ws.Cells.Item($r, 31).Value = 70
    $ws.Cells.Item($r, 32).Value = 0
}
